$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-04-21"

# Update the header label for the April row (shared string)
$ws.Range("A5").Value = "April (through 04-21)"

# Update the April row figures
$ws.Range("D5").Value = 41
$ws.Range("E5").Value = 38
$ws.Range("H5").Value = 78
$ws.Range("I5").Value = 91

# Update the Total row figures
$ws.Range("D6").Value = 230
$ws.Range("E6").Value = 235
$ws.Range("H6").Value = 501
$ws.Range("I6").Value = 527
